$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.817.95"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "2.499.36"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "588.08"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "176.53"
$ws.Range("E6").Value = "  +3.99%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("E9").Value = "  +3.92%  "
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "2.957.50"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "25.81"
$ws.Range("E14").Value = "  +2.10%  "
$ws.Range("D15").Value = "67.665.63"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").Value = "2.493.51"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "11.08"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "7.56"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").Value = "351.40"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "4.09"
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D23").Value = "70.88"
$ws.Range("E23").Value = "  +3.38%  "
$ws.Range("D24").Value = "4.32"
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("D25").Value = "1.79"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "9.18"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").Value = "2.626.45"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "0.0₃0911"
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("D30").Value = "509.30"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "7.84"
$ws.Range("E31").Value = "  +2.42%  "
$ws.Range("D32").Value = "1.26"
$ws.Range("E32").Value = "  +3.03%  "
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  +6.43%  "
$ws.Range("D36").Value = "161.98"
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("B37").Value = "WhiteBITCoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D37").Value = "18.68"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "18.42"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").Value = "1.34"
$ws.Range("E39").Value = "  +1.26%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  +3.50%  "
$ws.Range("D42").Value = "0.330"
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("D43").Value = "4.85"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("E44").Value = "  +2.80%  "
$ws.Range("D45").Value = "145.54"
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0256"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0743"
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").Value = "1.59"
$ws.Range("E50").Value = "  +1.78%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.588"
$ws.Range("E51").Value = "  +0.99%  "
